# "Revert to use older Op-Amp" - BOM row for D1's LED part reverted back to an
# older Digikey part (green 0603 LED instead of the UV flip-chip LED).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 (Ref-Des D1 / Function "Led") - swap in the older LED part info.
# Write the URL column first so new shared-string entries land in the same
# order Excel would append them (URL, Description, P/N, Digikey P/N).
$ws.Range("I8").Value = "http://www.digikey.com/product-search/en?keywords=160-1475-1-ND"
$ws.Range("C8").Value = "LED GREEN CLEAR 0603 R/A SMD"
$ws.Range("D8").Value = "LTST-S270GKT"
$ws.Range("E8").Value = "160-1475-1-ND"
$ws.Range("G8").Value = 0.1122

# The Extended Price column (H) already holds =G8*F8 as a shared formula,
# so it recalculates automatically once G8 changes.

# Reflect the author's last-selected cell on save.
$ws.Range("G8").Select()
